$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2 through 452) holds the "Förändrad" (last changed) date,
# stored as serial date 45178 (2023-09-09). Bump it by one day to 45179
# (2023-09-10) for every data row.
$ws.Range("C2:C452").Value = 45179
